$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the market report. It is inserted
# as a new row 10, pushing the previously existing rows 10-85 down to 11-86
# (row 85's data ends up, unchanged, in the new row 86).
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 44503
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 100112031
$ws.Cells.Item(10, 7).Value = "Poroto verde"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 65
$ws.Cells.Item(10, 11).Value = 43000
$ws.Cells.Item(10, 12).Value = 43000
$ws.Cells.Item(10, 13).Value = 43000
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 1720
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
